# Added API endpoints for receiving POST requests - backfilling newly
# received monthly data rows (2025-2028) in the itcz10 dataset sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 56 (2025): B56 already has a value; fill in C..M
$ws.Range("C56").Value = 3.27
$ws.Range("D56").Value = 4.68
$ws.Range("E56").Value = 8.91
$ws.Range("F56").Value = 10.01
$ws.Range("G56").Value = 10.93
$ws.Range("H56").Value = 15.16
$ws.Range("I56").Value = 13.86
$ws.Range("J56").Value = 9.59
$ws.Range("K56").Value = 9.56
$ws.Range("L56").Value = 7.94
$ws.Range("M56").Value = 2.53

# Row 57 (2026): fill in B..M
$ws.Range("B57").Value = 5.07
$ws.Range("C57").Value = 4.78
$ws.Range("D57").Value = 4.33
$ws.Range("E57").Value = 4.67
$ws.Range("F57").Value = 8.47
$ws.Range("G57").Value = 9.34
$ws.Range("H57").Value = 15.92
$ws.Range("I57").Value = 13.89
$ws.Range("J57").Value = 12.8
$ws.Range("K57").Value = 7.41
$ws.Range("L57").Value = 4.97
$ws.Range("M57").Value = 3.14

# Row 58 (2027): fill in B..M
$ws.Range("B58").Value = 0.95
$ws.Range("C58").Value = 2.07
$ws.Range("D58").Value = 3.11
$ws.Range("E58").Value = 7.39
$ws.Range("F58").Value = 7.98
$ws.Range("G58").Value = 11.72
$ws.Range("H58").Value = 10.63
$ws.Range("I58").Value = 10.5
$ws.Range("J58").Value = 12.65
$ws.Range("K58").Value = 8.85
$ws.Range("L58").Value = 7.6
$ws.Range("M58").Value = 3.23

# Row 59 (2028): fill in B only
$ws.Range("B59").Value = 6.21
